# TC07_SearchCategory.xlsx - "release part 5" update
#
# Two new "WAIT" keyword rows are inserted into the TC07_SearchCategory
# sheet's step table:
#   - one directly above the existing "CLICK / Gearing" row (old row 4)
#   - one directly above the existing "VERIFY_TEXT_PRESENT / GearingCategoryHeader" row
#     (old row 5, which has by then shifted down to row 6)
#
# Each inserted row only has a value in column B ("WAIT"); the rest of the
# row stays empty but keeps the same bordered formatting as the rest of the
# step table (copied down from the row directly above it, matching Excel's
# default "insert row" behaviour).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert first WAIT row (becomes row 4, pushes CLICK down to row 5) ---
$ws.Rows.Item(4).Insert()
$ws.Range("A3:E3").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)
$excel.CutCopyMode() = 0
$ws.Range("B4").Value() = "WAIT"

# --- Insert second WAIT row (becomes row 6, pushes VERIFY_TEXT_PRESENT down to row 7) ---
$ws.Rows.Item(6).Insert()
$ws.Range("A5:E5").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)
$excel.CutCopyMode() = 0
$ws.Range("B6").Value() = "WAIT"

# Leave the selection where the author ended up after making these edits.
$ws.Range("B11").Select() | Out-Null
